# "Fruta / hortaliza, semanal" - weekly refresh of the Piña (Feria Lagunitas
# de Puerto Montt) price series: a new week's record is inserted at the top
# of the data block (row 428), pushing the existing rows 428-448 down to
# 429-449 (the last former row, 448, ends up at 449).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 428 - this shifts rows 428:448 down to 429:449
# and extends the sheet's used range to row 449, exactly like Excel's
# Home > Insert > Insert Sheet Rows.
$ws.Rows.Item(428).Insert()

# Populate the new row with the latest week's observation.
$ws.Cells.Item(428, 1).Value  = 4
$ws.Cells.Item(428, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(428, 3).Value  = "Los Lagos"
$ws.Cells.Item(428, 4).Value  = 45147
$ws.Cells.Item(428, 5).Value  = 10
$ws.Cells.Item(428, 6).Value  = "Fruta"
$ws.Cells.Item(428, 7).Value  = 100108
$ws.Cells.Item(428, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(428, 9).Value  = 100108005
$ws.Cells.Item(428, 10).Value = "Piña"
$ws.Cells.Item(428, 11).Value = "Caramelo"
$ws.Cells.Item(428, 12).Value = "Segunda"
$ws.Cells.Item(428, 13).Value = 40
$ws.Cells.Item(428, 14).Value = 24000
$ws.Cells.Item(428, 15).Value = 24000
$ws.Cells.Item(428, 16).Value = 24000
$ws.Cells.Item(428, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(428, 18).Value = "Ecuador"
$ws.Cells.Item(428, 19).Value = 1714
$ws.Cells.Item(428, 20).Value = 14
